$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    # Force the cell to Text so numeric-looking strings (e.g. "13.20",
    # "4.00") keep their exact original formatting instead of being
    # parsed/normalized as numbers, then restore the default "Normal"
    # style so no stray formatting is left on the cell.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

Set-TextValue "D2" '62.224.30'
$ws.Range("E2").Value = '  +3.62%  '

Set-TextValue "D3" '3.415.62'
$ws.Range("E3").Value = '  +4.63%  '

$ws.Range("E4").Value = '  -0.03%  '

Set-TextValue "D5" '407.19'
$ws.Range("E5").Value = '  +0.38%  '

Set-TextValue "D6" '131.85'
$ws.Range("E6").Value = '  +19.11%  '

Set-TextValue "D7" '0.611'
$ws.Range("E7").Value = '  +8.84%  '

$ws.Range("E8").Value = '  +0.00%  '

Set-TextValue "D9" '0.679'
$ws.Range("E9").Value = '  +10.86%  '

Set-TextValue "D10" '0.128'
$ws.Range("E10").Value = '  +13.97%  '

Set-TextValue "D11" '42.39'
$ws.Range("E11").Value = '  +10.67%  '

$ws.Range("E12").Value = '  -0.24%  '

Set-TextValue "D13" '3.960.39'
$ws.Range("E13").Value = '  +4.58%  '

Set-TextValue "D14" '8.59'
$ws.Range("E14").Value = '  +6.49%  '

$ws.Range("E15").Value = '  +5.11%  '

Set-TextValue "D16" '3.399.50'
$ws.Range("E16").Value = '  +2.01%  '

Set-TextValue "D17" '62.070.84'
$ws.Range("E17").Value = '  +3.30%  '

Set-TextValue "D18" '11.54'
$ws.Range("E18").Value = '  +10.18%  '

$ws.Range("E19").Value = '  +5.86%  '

$ws.Range("E20").Value = '  +19.53%  '

Set-TextValue "D21" '3.28'
$ws.Range("E21").Value = '  +0.89%  '

Set-TextValue "D22" '82.98'
$ws.Range("E22").Value = '  +13.95%  '

Set-TextValue "D23" '13.20'
$ws.Range("E23").Value = '  +7.15%  '

Set-TextValue "D24" '309.11'
$ws.Range("E24").Value = '  +5.22%  '

$ws.Range("E25").Value = '  +3.72%  '

Set-TextValue "D26" '8.58'
$ws.Range("E26").Value = '  +16.04%  '

Set-TextValue "D27" '29.85'
$ws.Range("E27").Value = '  +4.08%  '

Set-TextValue "D28" '4.64'
$ws.Range("E28").Value = '  +8.94%  '

$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue "D29" '0.175'
$ws.Range("E29").Value = '  +2.12%  '

$ws.Range("B30").Value = 'RenderToken'
$ws.Range("C30").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D30" '7.49'
$ws.Range("E30").Value = '  +2.08%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
Set-TextValue "D31" '0.116'
$ws.Range("E31").Value = '  +3.31%  '

$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue "D32" '11.86'
$ws.Range("E32").Value = '  +7.26%  '

$ws.Range("E33").Value = '  +7.90%  '

Set-TextValue "D34" '42.90'
$ws.Range("E34").Value = '  +11.37%  '

$ws.Range("E35").Value = '  +0.07%  '

Set-TextValue "D36" '0.0490'
$ws.Range("E36").Value = '  +3.31%  '

Set-TextValue "D37" '52.44'
$ws.Range("E37").Value = '  +0.91%  '

Set-TextValue "D38" '0.996'
$ws.Range("E38").Value = '  -0.43%  '

Set-TextValue "D39" '3.45'
$ws.Range("E39").Value = '  +5.53%  '

$ws.Range("E40").Value = '  -2.54%  '

$ws.Range("E41").Value = '  +10.59%  '

$ws.Range("E42").Value = '  +6.26%  '

Set-TextValue "D43" '137.55'
$ws.Range("E43").Value = '  +1.88%  '

Set-TextValue "D44" '4.00'
$ws.Range("E44").Value = '  +7.60%  '

$ws.Range("B45").Value = 'Celestia'
$ws.Range("C45").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-TextValue "D45" '17.16'
$ws.Range("E45").Value = '  +6.90%  '

$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D46" '0.287'
$ws.Range("E46").Value = '  -0.58%  '

Set-TextValue "D47" '2.23'
$ws.Range("E47").Value = '  +2.01%  '

Set-TextValue "D48" '21.80'
$ws.Range("E48").Value = '  +5.41%  '

Set-TextValue "D49" '2.153.55'
$ws.Range("E49").Value = '  +2.40%  '

Set-TextValue "D50" '3.747.86'
$ws.Range("E50").Value = '  +4.02%  '

$ws.Range("E51").Value = '  -0.68%  '
